$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (row) was added to the log. Excel keeps the sheet
# sorted with the newest entries near the top of the data block, so the
# existing rows 3:16 get pushed down one row and the fresh record lands
# in the now-empty row 3.
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44532
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 300000000
$ws.Cells.Item(3, 7).Value = "Espárragos"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 240
$ws.Cells.Item(3, 11).Value = 800
$ws.Cells.Item(3, 12).Value = 900
$ws.Cells.Item(3, 13).Value = 850
$ws.Cells.Item(3, 14).Value = "`$/kilo"
$ws.Cells.Item(3, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(3, 16).Value = 850
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# D3 carries the same date-style formatting as the other Fecha cells.
$ws.Cells.Item(3, 4).Style = $ws.Cells.Item(4, 4).Style
